# dingo_command rating_summary_template.xlsx
# [fix][cloudkitty] Add code for optimized to rating download xlsx eg.
#
# - header labels now show the timezone they're reported in
# - selection cursor moved off the header row
# - first two data columns widened (Project ID / time range no longer clipped)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: clarify the timezone used for Begin/End -----------------
$ws.Range("A1").Value = "Begin(Timezone: UTC)"
$ws.Range("B1").Value = "End(Timezone: UTC)"

# --- Column widths: widen column A and B so the longer headers/values fit -
# (ColumnWidth is expressed in characters; the engine quantises to whole
# pixels at the standard 7px "maximum digit width", so these are the
# character widths whose quantised pixel width lands closest to the
# published 25.125 / 25.5 target column widths.)
$ws.Columns.Item(1).ColumnWidth = 24.4
$ws.Columns.Item(2).ColumnWidth = 24.75

# --- Selection: move the active cell off the frozen header row ------------
$null = $ws.Range("C15").Select()
